$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -1.63190989412726
$ws.Cells.Item(2, 3).Value = 0.376351439199373
$ws.Cells.Item(2, 4).Value = -4.33613299738906
$ws.Cells.Item(2, 5).Value = 0.0000145011223854152
$ws.Cells.Item(2, 6).Value = 0.004280958796761

$ws.Cells.Item(3, 2).Value = -1.35634275007098
$ws.Cells.Item(3, 3).Value = 0.164297922019916
$ws.Cells.Item(3, 4).Value = -8.25538590747705
$ws.Cells.Item(3, 5).Value = 0.000000000000000151412204669717
$ws.Cells.Item(3, 6).Value = 0.000000000000759887384502422

$ws.Cells.Item(4, 2).Value = -2.94679643856939
$ws.Cells.Item(4, 3).Value = 0.440101590221266
$ws.Cells.Item(4, 4).Value = -6.69571867960725
$ws.Cells.Item(4, 5).Value = 0.000000000021461424338273
$ws.Cells.Item(4, 6).Value = 0.0000000461604578338626

$ws.Cells.Item(5, 2).Value = -1.62032219173922
$ws.Cells.Item(5, 3).Value = 0.348781598260943
$ws.Cells.Item(5, 4).Value = -4.64566422029803
$ws.Cells.Item(5, 5).Value = 0.00000338983988562143
$ws.Cells.Item(5, 6).Value = 0.00150110086229165

$ws.Cells.Item(6, 2).Value = -1.24134152870974
$ws.Cells.Item(6, 3).Value = 0.251082177019885
$ws.Cells.Item(6, 4).Value = -4.94396513302269
$ws.Cells.Item(6, 5).Value = 0.00000076549400353011
$ws.Cells.Item(6, 6).Value = 0.000411617061326762

$ws.Cells.Item(7, 2).Value = 1.588823223888
$ws.Cells.Item(7, 3).Value = 0.348761029744535
$ws.Cells.Item(7, 4).Value = 4.55562143812858
$ws.Cells.Item(7, 5).Value = 0.00000522309293635477
$ws.Cells.Item(7, 6).Value = 0.00215596806320689

$ws.Cells.Item(8, 2).Value = -1.90772500328651
$ws.Cells.Item(8, 3).Value = 0.439109029584707
$ws.Cells.Item(8, 4).Value = -4.3445360371905
$ws.Cells.Item(8, 5).Value = 0.0000139570331069236
$ws.Cells.Item(8, 6).Value = 0.00420274180915683

$ws.Cells.Item(9, 2).Value = -1.44823364645923
$ws.Cells.Item(9, 3).Value = 0.241582286064411
$ws.Cells.Item(9, 4).Value = -5.99478409635177
$ws.Cells.Item(9, 5).Value = 0.00000000203755965260201
$ws.Cells.Item(9, 6).Value = 0.00000235980754842892

$ws.Cells.Item(10, 2).Value = -1.60356836088316
$ws.Cells.Item(10, 3).Value = 0.214974356518554
$ws.Cells.Item(10, 4).Value = -7.4593471837873
$ws.Cells.Item(10, 5).Value = 0.00000000000008695221192634
$ws.Cells.Item(10, 6).Value = 0.000000000218192083793829

$ws.Cells.Item(11, 2).Value = -1.01348281889532
$ws.Cells.Item(11, 3).Value = 0.258127047586793
$ws.Cells.Item(11, 4).Value = -3.92629454514852
$ws.Cells.Item(11, 5).Value = 0.0000862645105444369
$ws.Cells.Item(11, 6).Value = 0.0170894535625927

$ws.Cells.Item(12, 2).Value = -1.04147988658406
$ws.Cells.Item(12, 3).Value = 0.180639932696432
$ws.Cells.Item(12, 4).Value = -5.76550196314724
$ws.Cells.Item(12, 5).Value = 0.00000000814151543858835
$ws.Cells.Item(12, 6).Value = 0.00000766116602771164

$ws.Cells.Item(13, 2).Value = -1.33955827239898
$ws.Cells.Item(13, 3).Value = 0.299172832705637
$ws.Cells.Item(13, 4).Value = -4.47753982299926
$ws.Cells.Item(13, 5).Value = 0.00000755081455570635
$ws.Cells.Item(13, 6).Value = 0.00277280643782231

$ws.Cells.Item(14, 2).Value = -1.18530194099526
$ws.Cells.Item(14, 3).Value = 0.292784479057877
$ws.Cells.Item(14, 4).Value = -4.04837696591475
$ws.Cells.Item(14, 5).Value = 0.0000515740332346349
$ws.Cells.Item(14, 6).Value = 0.0117651309754646

$ws.Cells.Item(15, 2).Value = -1.07909394962812
$ws.Cells.Item(15, 3).Value = 0.203999656550554
$ws.Cells.Item(15, 4).Value = -5.28968512925271
$ws.Cells.Item(15, 5).Value = 0.000000122527108209975
$ws.Cells.Item(15, 6).Value = 0.0000922384070604693

$ws.Cells.Item(16, 2).Value = 1.07831751647532
$ws.Cells.Item(16, 3).Value = 0.177194952408177
$ws.Cells.Item(16, 4).Value = 6.08548664519159
$ws.Cells.Item(16, 5).Value = 0.00000000116138005404124
$ws.Cells.Item(16, 6).Value = 0.00000158961255396772

$ws.Cells.Item(17, 2).Value = -1.15268260947424
$ws.Cells.Item(17, 3).Value = 0.234176510717812
$ws.Cells.Item(17, 4).Value = -4.92228108592519
$ws.Cells.Item(17, 5).Value = 0.000000855412292949811
$ws.Cells.Item(17, 6).Value = 0.000429302916088412

$ws.Cells.Item(18, 2).Value = 1.27464828373868
$ws.Cells.Item(18, 3).Value = 0.338214173997842
$ws.Cells.Item(18, 4).Value = 3.76876068992549
$ws.Cells.Item(18, 5).Value = 0.000164060083868864
$ws.Cells.Item(18, 6).Value = 0.0265600927175228

$ws.Cells.Item(19, 2).Value = -2.03364265983309
$ws.Cells.Item(19, 3).Value = 0.527117562945501
$ws.Cells.Item(19, 4).Value = -3.85804382701501
$ws.Cells.Item(19, 5).Value = 0.000114298155773964
$ws.Cells.Item(19, 6).Value = 0.0204865837301523

$ws.Cells.Item(20, 2).Value = -1.01231252706982
$ws.Cells.Item(20, 3).Value = 0.23423772199123
$ws.Cells.Item(20, 4).Value = -4.3217314378925
$ws.Cells.Item(20, 5).Value = 0.0000154809551702853
$ws.Cells.Item(20, 6).Value = 0.00432848243844727

$ws.Cells.Item(21, 2).Value = -1.2998362358837
$ws.Cells.Item(21, 3).Value = 0.292149183042641
$ws.Cells.Item(21, 4).Value = -4.44922084787751
$ws.Cells.Item(21, 5).Value = 0.00000861823630871352
$ws.Cells.Item(21, 6).Value = 0.00308943252057121

$ws.Cells.Item(22, 2).Value = 7.49574097361223
$ws.Cells.Item(22, 3).Value = 1.64646870221189
$ws.Cells.Item(22, 4).Value = 4.55261673880734
$ws.Cells.Item(22, 5).Value = 0.00000529827433173849
$ws.Cells.Item(22, 6).Value = 0.00215596806320689

$ws.Cells.Item(23, 2).Value = -1.28503409099958
$ws.Cells.Item(23, 3).Value = 0.327267527418009
$ws.Cells.Item(23, 4).Value = -3.92655544269214
$ws.Cells.Item(23, 5).Value = 0.0000861710363757061
$ws.Cells.Item(23, 6).Value = 0.0170894535625927

$ws.Cells.Item(24, 2).Value = -1.37568665210752
$ws.Cells.Item(24, 3).Value = 0.37408012598569
$ws.Cells.Item(24, 4).Value = -3.67751868261547
$ws.Cells.Item(24, 5).Value = 0.000235513830914385
$ws.Cells.Item(24, 6).Value = 0.034426177070359

$ws.Cells.Item(25, 2).Value = 2.62858998443568
$ws.Cells.Item(25, 3).Value = 0.464741305694865
$ws.Cells.Item(25, 4).Value = 5.6560283155927
$ws.Cells.Item(25, 5).Value = 0.0000000154915920778827
$ws.Cells.Item(25, 6).Value = 0.0000137200829602707

$ws.Cells.Item(26, 2).Value = 2.13947173668205
$ws.Cells.Item(26, 3).Value = 0.268332388978483
$ws.Cells.Item(26, 4).Value = 7.97321465674279
$ws.Cells.Item(26, 5).Value = 0.00000000000000154599307390871
$ws.Cells.Item(26, 6).Value = 0.00000000000465529434415389

$ws.Cells.Item(27, 2).Value = 1.37457298343047
$ws.Cells.Item(27, 3).Value = 0.219706399277337
$ws.Cells.Item(27, 4).Value = 6.25640849766664
$ws.Cells.Item(27, 5).Value = 0.000000000393944061791139
$ws.Cells.Item(27, 6).Value = 0.000000659024643814154

$ws.Cells.Item(28, 2).Value = -1.4300251279996
$ws.Cells.Item(28, 3).Value = 0.289697245546104
$ws.Cells.Item(28, 4).Value = -4.93627450721487
$ws.Cells.Item(28, 5).Value = 0.000000796290074761562
$ws.Cells.Item(28, 6).Value = 0.000413411840193451

$ws.Cells.Item(29, 2).Value = 1.19130997794583
$ws.Cells.Item(29, 3).Value = 0.23253678887264
$ws.Cells.Item(29, 4).Value = 5.12310324624938
$ws.Cells.Item(29, 5).Value = 0.000000300547668919883
$ws.Cells.Item(29, 6).Value = 0.000215478366821798

$ws.Cells.Item(30, 2).Value = -1.02680552312415
$ws.Cells.Item(30, 3).Value = 0.20426241548003
$ws.Cells.Item(30, 4).Value = -5.02689406032477
$ws.Cells.Item(30, 5).Value = 0.000000498487551529602
$ws.Cells.Item(30, 6).Value = 0.000312717857326237

$ws.Cells.Item(31, 2).Value = -1.15717732423933
$ws.Cells.Item(31, 3).Value = 0.254613362874554
$ws.Cells.Item(31, 4).Value = -4.54484128866979
$ws.Cells.Item(31, 5).Value = 0.00000549766625253847
$ws.Cells.Item(31, 6).Value = 0.00217823323942682

$ws.Cells.Item(32, 2).Value = -1.02392661796103
$ws.Cells.Item(32, 3).Value = 0.231649245933928
$ws.Cells.Item(32, 4).Value = -4.42015951242542
$ws.Cells.Item(32, 5).Value = 0.00000986280719515315
$ws.Cells.Item(32, 6).Value = 0.00345335872395874

$ws.Cells.Item(33, 2).Value = 1.02182974976917
$ws.Cells.Item(33, 3).Value = 0.236695711984775
$ws.Cells.Item(33, 4).Value = 4.31706067338854
$ws.Cells.Item(33, 5).Value = 0.0000158120705442747
$ws.Cells.Item(33, 6).Value = 0.00432848243844727

$ws.Cells.Item(34, 2).Value = -1.13481150507729
$ws.Cells.Item(34, 3).Value = 0.320260168941147
$ws.Cells.Item(34, 4).Value = -3.54340506604127
$ws.Cells.Item(34, 5).Value = 0.00039499549400244
$ws.Cells.Item(34, 6).Value = 0.0479600980459737

$ws.Cells.Item(35, 2).Value = -1.08234934698243
$ws.Cells.Item(35, 3).Value = 0.218860608405867
$ws.Cells.Item(35, 4).Value = -4.94538215381023
$ws.Cells.Item(35, 5).Value = 0.000000759946327659539
$ws.Cells.Item(35, 6).Value = 0.000411617061326762

$ws.Cells.Item(36, 2).Value = 1.21333749240025
$ws.Cells.Item(36, 3).Value = 0.201767101641471
$ws.Cells.Item(36, 4).Value = 6.01355465053112
$ws.Cells.Item(36, 5).Value = 0.00000000181498703506128
$ws.Cells.Item(36, 6).Value = 0.00000227720373332355

$ws.Cells.Item(37, 2).Value = 1.16866816731873
$ws.Cells.Item(37, 3).Value = 0.189511167334993
$ws.Cells.Item(37, 4).Value = 6.16675092952655
$ws.Cells.Item(37, 5).Value = 0.000000000697074246686327
$ws.Cells.Item(37, 6).Value = 0.00000104951498581093

$ws.Cells.Item(38, 2).Value = -1.54100616336134
$ws.Cells.Item(38, 3).Value = 0.341143674949919
$ws.Cells.Item(38, 4).Value = -4.51717641720183
$ws.Cells.Item(38, 5).Value = 0.0000062669691077204
$ws.Cells.Item(38, 6).Value = 0.00241937145861124

$ws.Cells.Item(39, 2).Value = -1.46229451754072
$ws.Cells.Item(39, 3).Value = 0.401853199185281
$ws.Cells.Item(39, 4).Value = -3.63887738235102
$ws.Cells.Item(39, 5).Value = 0.0002738291331149
$ws.Cells.Item(39, 6).Value = 0.0378235910842013

$ws.Cells.Item(40, 2).Value = -1.81067359782906
$ws.Cells.Item(40, 3).Value = 0.465765878465489
$ws.Cells.Item(40, 4).Value = -3.8875187761596
$ws.Cells.Item(40, 5).Value = 0.000101274155950972
$ws.Cells.Item(40, 6).Value = 0.0193010593923776

$ws.Cells.Item(41, 2).Value = -1.24409954428466
$ws.Cells.Item(41, 3).Value = 0.344830790169355
$ws.Cells.Item(41, 4).Value = -3.60785515607134
$ws.Cells.Item(41, 5).Value = 0.000308738757863268
$ws.Cells.Item(41, 6).Value = 0.0407800402751629
